$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 430
$ws.Range("F5").Value = 3950
$ws.Range("F7").Value = 2603
$ws.Range("F11").Value = 2346
$ws.Range("F14").Value = 331
$ws.Range("F16").Value = 24
$ws.Range("F19").Value = 359
$ws.Range("F20").Value = 317
$ws.Range("F21").Value = 468
$ws.Range("F22").Value = 687
$ws.Range("F24").Value = 437
$ws.Range("F25").Value = 15
$ws.Range("F26").Value = 1313
$ws.Range("F28").Value = 175
$ws.Range("F29").Value = 10
$ws.Range("F31").Value = 4461
$ws.Range("F32").Value = 4363
$ws.Range("F34").Value = 309
$ws.Range("F35").Value = 78
$ws.Range("F39").Value = 13
$ws.Range("F40").Value = 510
$ws.Range("F43").Value = 190
$ws.Range("F47").Value = 72

$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 27
$ws.Range("F14").Value = 46
$ws.Range("F15").Value = 219

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 2346
$ws.Range("F5").Value = 59

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 430
$ws.Range("F6").Value = 3950
$ws.Range("F8").Value = 2603
$ws.Range("F12").Value = 2346
$ws.Range("F15").Value = 331
$ws.Range("F17").Value = 24
$ws.Range("F20").Value = 359
$ws.Range("F21").Value = 687
$ws.Range("F23").Value = 1313
$ws.Range("F27").Value = 27
$ws.Range("F29").Value = 4461
$ws.Range("F30").Value = 4363
$ws.Range("F35").Value = 13
$ws.Range("F38").Value = 510
$ws.Range("F42").Value = 46
$ws.Range("F44").Value = 190
$ws.Range("F47").Value = 72
$ws.Range("F49").Value = 219
